$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Bottom of the document: the paragraph that duplicated the H1 title
#    ("Play Action Jack Slot for Free ...") is removed entirely, and the
#    following (italic) paragraph's text is replaced with the new
#    feature-image prompt, while keeping its italic formatting.
# -----------------------------------------------------------------------

$dupTitlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.Contains("Play Action Jack Slot for Free") -and $t.Contains("Fun Jungle Theme")) {
        $dupTitlePara = $p
    }
}

if ($dupTitlePara -ne $null) {
    $dupTitlePara.Range.Delete()
}

# Re-find the italic paragraph (its index shifted after the delete above,
# and its Range object may be stale), then swap its text in place so the
# existing run formatting (italic) is preserved.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains("Read our review of Action Jack slot machine")) {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        $newPrompt = 'Create an engaging feature image for "Action Jack" that captures the excitement and adventure of the game. The image should be in cartoon style and feature a happy Maya warrior with glasses, as this is the main character of the game. The Maya warrior should be shown in the depths of the jungle, surrounded by symbols from the slot game such as treasure maps, lava waterfalls, gates, amulets, mysterious statues, and erupting volcanoes. The Maya warrior should be depicted as confident and fearless, ready to take on any challenge in his quest for ancient treasures. The image should be eye-catching and colorful, with the Maya warrior prominently displayed in the center of the graphic. '
        $target = $d.Range($pStart, $pEnd - 1)
        $target.Text = $newPrompt
        break
    }
}

# -----------------------------------------------------------------------
# 2) Top of the document: insert a new "Meta description" paragraph right
#    after the Heading1 title paragraph.
# -----------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Range.ParagraphFormat.Style = "Normal"

$metaStart = $metaPara.Range.Start
$metaEndInit = $metaPara.Range.End
$metaFullText = "Meta description: Read our review of Action Jack slot machine, play for free, and discover its low volatility, fun design, and Free Spins Bonus with 3x multiplier."
$metaTarget = $d.Range($metaStart, $metaEndInit - 1)
$metaTarget.Text = $metaFullText

$metaBoldRange = $d.Range($metaStart, $metaStart + 16)
$metaBoldRange.Bold = 1

Write-Output "done"
